$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.436.39'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '3.168.59'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.165.32'
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.448'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.73%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("D13").Value = '3.718.29'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000167'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.05%  '
$ws.Range("D17").Value = '58.517.32'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '3.168.68'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '360.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.22%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("D27").Value = '0.0₃0959'
$ws.Range("E27").Value = '  +5.49%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  +4.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  +2.59%  '
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("E40").Value = '  +13.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0677'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("E42").Value = '  +4.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.95%  '
$ws.Range("D44").Value = '3.212.03'
$ws.Range("E44").Value = '  +2.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0274'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '36.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = '2.351.62'
$ws.Range("E47").Value = '  +2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  +5.90%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  +1.48%  '
